$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "26.716.06"
$ws.Range('E2').Value = "  +0.96%  "
$ws.Range('D3').Value = "1.647.75"
$ws.Range('E3').Value = "  +1.25%  "
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = "  +0.05%  "
$ws.Range('D5').Value = "'216.11"
$ws.Range('E5').Value = "  +1.36%  "
$ws.Range('D6').Value = "'0.507"
$ws.Range('E6').Value = "  +1.59%  "
$ws.Range('E7').Value = "  +0.35%  "
$ws.Range('E8').Value = "  +1.28%  "
$ws.Range('D9').Value = "'0.0627"
$ws.Range('E9').Value = "  +0.48%  "
$ws.Range('D10').Value = "'19.15"
$ws.Range('E10').Value = "  +1.25%  "
$ws.Range('D11').Value = "'0.0844"
$ws.Range('E11').Value = "  +0.04%  "
$ws.Range('D12').Value = "1.877.26"
$ws.Range('E12').Value = "  +1.27%  "
$ws.Range('D13').Value = "1.668.52"
$ws.Range('E13').Value = "  +1.58%  "
$ws.Range('D14').Value = "'4.18"
$ws.Range('E14').Value = "  +1.16%  "
$ws.Range('D15').Value = "'0.532"
$ws.Range('E15').Value = "  +1.82%  "
$ws.Range('D16').Value = "'65.14"
$ws.Range('E16').Value = "  +0.51%  "
$ws.Range('D17').Value = "26.717.55"
$ws.Range('E17').Value = "  +0.83%  "
$ws.Range('D18').Value = "0.0₃0743"
$ws.Range('E18').Value = "  +0.38%  "
$ws.Range('D19').Value = "'218.41"
$ws.Range('E19').Value = "  +1.60%  "
$ws.Range('E20').Value = "  +0.38%  "
$ws.Range('D21').Value = "'4.36"
$ws.Range('E21').Value = "  +1.34%  "
$ws.Range('D22').Value = "'6.28"
$ws.Range('E22').Value = "  +0.24%  "
$ws.Range('B23').Value = "Avalanche"
$ws.Range('C23').Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('D23').Value = "'9.50"
$ws.Range('E23').Value = "  +2.43%  "
$ws.Range('B24').Value = "Toncoin"
$ws.Range('C24').Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('D24').Value = "'2.27"
$ws.Range('E24').Value = "  +14.04%  "
$ws.Range('D25').Value = "'145.99"
$ws.Range('E25').Value = "  -1.69%  "
$ws.Range('E26').Value = "  +0.32%  "
$ws.Range('E27').Value = "  +0.17%  "
$ws.Range('D28').Value = "'7.13"
$ws.Range('E28').Value = "  +4.33%  "
$ws.Range('E29').Value = "  +0.99%  "
$ws.Range('D30').Value = "'0.0515"
$ws.Range('E30').Value = "  +1.32%  "
$ws.Range('D31').Value = "'1.18"
$ws.Range('E31').Value = "  +1.51%  "
$ws.Range('D32').Value = "'3.37"
$ws.Range('E32').Value = "  +1.42%  "
$ws.Range('D33').Value = "'3.01"
$ws.Range('E33').Value = "  +2.16%  "
$ws.Range('D34').Value = "1.277.47"
$ws.Range('E34').Value = "  +4.77%  "
$ws.Range('E35').Value = "  +3.54%  "
$ws.Range('D36').Value = "'2.43"
$ws.Range('E36').Value = "  +1.75%  "
$ws.Range('D37').Value = "'0.0179"
$ws.Range('E37').Value = "  +3.03%  "
$ws.Range('D38').Value = "'0.534"
$ws.Range('E38').Value = "  +5.77%  "
$ws.Range('D39').Value = "'0.823"
$ws.Range('E39').Value = "  +3.45%  "
$ws.Range('E40').Value = "  +0.43%  "
$ws.Range('D41').Value = "'0.815"
$ws.Range('E41').Value = "  +2.65%  "
$ws.Range('D42').Value = "'2.25"
$ws.Range('E42').Value = "  +0.02%  "
$ws.Range('D43').Value = "'5.45"
$ws.Range('E43').Value = "  +1.69%  "
$ws.Range('D44').Value = "1.785.99"
$ws.Range('D45').Value = "'91.85"
$ws.Range('E45').Value = "  -1.40%  "
$ws.Range('D46').Value = "'59.86"
$ws.Range('E47').Value = "  +1.21%  "
$ws.Range('E48').Value = "  +1.11%  "
$ws.Range('D49').Value = "'7.75"
$ws.Range('E49').Value = "  +3.07%  "
$ws.Range('D50').Value = "'0.0969"
$ws.Range('E50').Value = "  +1.98%  "
$ws.Range('D51').Value = "'0.407"
$ws.Range('E51').Value = "  +0.13%  "
